$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('https://trsym.com', 'info@tavarone.com')
    ,@('https://gypnqn.com.ar', 'info@gypnqn.com.ar')
    ,@('https://mg54.com', 'info@mg54.com')
    ,@('https://achavalcornejo.com', 'info@achavalcornejo.com')
    ,@('https://medihome.com.ar', 'info@medihome.com.ar')
    ,@('https://gizmosite.tv', 'jobs@gizmosite.tv')
    ,@('https://etman.com.ar', 'info@estacionvisual.com.ar')
    ,@('https://p-per.com', 'info@p-per.com')
    ,@('https://midios.com', 'hello@tyler.com, info@midios.com')
    ,@('https://friolatina.com', 'info@grupoltn.com, INFO@GRUPOLTN.COM')
    ,@('https://gruposolanas.com', 'marketing@gruposolanas.com')
    ,@('https://petropack.com', 'info@petropack.com')
    ,@('https://ase.com.ar', 'info@opticaconecta.com.ar')
    ,@('https://monitorclaves.com.ar', 'info@claves.com.ar')
    ,@('https://royalchef.com.ar', 'info@royalchef.com.ar')
    ,@('https://gruposbs.com', 'info@gruposbs.com')
    ,@('https://tiarg.com', 'info@tiarg.com.ar')
    ,@('https://iprofesional.com', 'info@iprofesional.com')
    ,@('https://bluetrailsoft.com', 'sales@bluetrail.software')
    ,@('https://linceseguridad.com.ar', 'info@linceseguridad.com.ar')
    ,@('https://sintectur.com.ar', 'info@sintectur.com.ar')
    ,@('https://englishservices.com.ar', 'info@englishservices.com.ar')
    ,@('https://emtech.com.ar', 'info@emtech.com.ar')
    ,@('https://mastersoft.com.ar', 'info@mastersoft.com.ar, Info@mastersoft.com.ar')
    ,@('https://zetti.com.ar', 'info@zetti.com.ar')
    ,@('https://lylyk.com.ar', 'info@lylyk.com.ar')
    ,@('https://ladolfina.com', 'info@ladolfina.com')
    ,@('https://elinpar.com', 'info@elinpar.com')
    ,@('https://riera.com.ar', 'info@riera.com.ar')
    ,@('https://colinet.com.ar', 'info@colinet.com.ar')
    ,@('https://istorming.com.ar', 'info@istorming.com.ar')
    ,@('https://corpomedica.com.ar', 'info@corpomedica.com.ar')
    ,@('https://we-are-around.com', 'info@we-are-around.com')
    ,@('https://ngsoluciones.com', 'info@ngsoluciones.com')
    ,@('https://ccmatienzo.com.ar', 'info@ccmatienzo.com.ar, info@matienzo.com.ar')
    ,@('https://exisoft.com.ar', 'info@exisoft.com.ar')
    ,@('https://feyalegria.org.ar', 'info@feyalegria.org.ar')
    ,@('https://solinf-sa.com.ar', 'info@solinf-sa.com.ar')
    ,@('https://ecopreneur.com.ar', 'info@ecopreneur.com.uy, info@ecopreneur.com.py, info@ecopreneur.com.bo, info@ecopreneur.com.ar')
    ,@('https://agidea.com.ar', 'info@agidea.com.ar')
    ,@('https://sooft.com.ar', 'info@sooft.tech')
    ,@('https://bfu.com.ar', 'info@bfu.com.ar')
    ,@('https://sanatoriodelsalvador.com', 'info@sanatoriodelsalvador.com')
    ,@('https://plenty.tv', 'hello@plentymore.agency')
    ,@('https://wetcorporate.com', 'info@wetcorporate.com, INFO@WETCORPORATE.COM')
    ,@('https://grupoeconat.com', 'info@grupoeconat.com')
    ,@('https://biwares.com', 'info@biwares.com')
    ,@('https://rrpp.org.ar', 'info@feedbackpr.com.ar, info@identiapr.com, Info@rpmcomunicaciones.com.ar, marketing@pulpolatam.com, info@grupomass.net, info@speconsultores.com.ar, info@ketchum.com.ar, info@zelmirakconsultora.com, info@arenapublicaconsultores.com, info@fullpr.com.ar, info@lanusse.ar, info@ballero.com.ar, info@moscuagencia.com, info@rpmcomunicaciones.com.ar, info@capitolpac.com, info@koifman.com.ar, info@abcomunicaciones.com, info@motivacom.com, info@vgv.com.ar')
    ,@('https://maiten.com', 'info@maiten.com')
    ,@('https://bolsasanjuan.com.ar', 'info@bolsasanjuan.com.ar')
    ,@('https://emaservicios.com.ar', 'info@emaservicios.com.ar')
    ,@('https://furlongincoming.com.ar', 'marketing@furlongincoming.com.ar, sales@furlongincoming.com.uy, sales@furlongincoming.com.ar')
    ,@('https://silicanetworks.com', 'info@silicanetworks.com')
    ,@('https://nexo-group.com.ar', 'Info@nexo-group.com.ar')
    ,@('https://evolutionp.com', 'info@evolutionp.com, INFO@EVOLUTIONP.COM')
    ,@('https://isf-argentina.org', 'info@isf-argentina.org')
    ,@('https://pinturasvenier.com', 'info@pinturasvenier.com.ar')
    ,@('https://assekuransa.com', 'office@mancas.ro, info@assekuransa.com, info@seaman-jo.com, info@aet-lloydsagency.com, info@ba-ty.com, info@mtscairo.com, info@mtsalexandria.com, info@mjkaquamarine.ie, info@bkk.lt, office@transpack.co.at, office@theokoch.com, info@sadaoui.net, info@castellano1869.it, info@demicheli-wassmuth.it, info@marinesurveyors.se, info@werner-brusberg-sl.com, info@liquidadores.net, info@interlloydclaims.nl, info@bdv.be, info@bianchi.com.mt, info@africamarine.be, info@cargoinspect.com, info@ostseekontor.de, info@mcleans.fr, admin@rsgamble.com, info@loadjusters.com, info@mtsportsaid.com, info@igis-surveyor.com, info@wesmans.com, info@assekeurope.es, office@scua.co.jp, office@vitsanbaku.com')
    ,@('https://cessi.org.ar', 'info@cessi.org.ar, info@prodasoftware.com')
    ,@('https://colegiosangabriel.com', 'info@centrodeportivosangabriel.es')
    ,@('https://room23.com.ar', 'hello@room23.com.ar')
    ,@('https://espaciobuenosaires.com.ar', 'info@espaciobuenosaires.com.ar, info@espaciobuenosaires.com')
    ,@('https://telectronica.com', 'info@telectronica.com')
    ,@('https://grobocopatel.com.ar', 'info@grobocopatel.com.ar')
    ,@('https://unobahiaclub.com', 'info@unobahiaclub.com')
    ,@('https://pehuenconsulting.com', 'info@pehuen.com')
    ,@('https://aedin.org', 'info@aedin.org')
    ,@('https://vidasilvestre.org.ar', 'info@vidasilvestre.org.ar')
    ,@('https://jpof.com.ar', 'info@jpof.com.ar')
    ,@('https://agrositio.com', 'info@agrositio.com, info@edn.org.ar')
    ,@('https://qservices.com.ar', 'info@qservices.com.ar')
    ,@('https://milhousehostel.com', 'info@milhousehostel.com')
    ,@('https://gemco.com.ar', 'info@gemco.com.ar')
    ,@('https://ecolatina.com', 'info@ecolatina.com')
    ,@('https://marketec.com.ar', 'info@marketec.com.ar, info@smark.com.mx')
    ,@('https://institucionulloa.com.ar', 'info@institucionulloa.com.ar')
    ,@('https://medoro.com.ar', 'info@website.com')
    ,@('https://qmamericas.com', 'info@qmamericas.com')
    ,@('https://destefano.com.ar', 'marketing@DESTEFANO1913.COM')
    ,@('https://ensalud.org', 'info@ensalud.org')
    ,@('https://other-guys.com', 'support@other-guys.com, sales@other-guys.com')
    ,@('https://fundaciondiscar.org.ar', 'info@fundaciondiscar.org.ar')
    ,@('https://artekstands.com', 'info@artekstands.com')
    ,@('https://entercommla.com', 'info@entercommla.com')
    ,@('https://cima-it.com', 'info@cima-it.com, support@cima-it.com')
    ,@('https://cervantes.edu.ar', 'info@cervantes.edu.ar')
    ,@('https://laboratoriosima.com.ar', 'info@laboratoriosima.com.ar')
    ,@('https://bluedraft.com.ar', 'info@example.com')
    ,@('https://soicos.com', 'info@soicos.com, billing@soicos.com')
    ,@('https://garridolawfirm.com', 'info@garridolawfirm.com')
    ,@('https://adblickagro.com', 'info@adblickagro.com')
    ,@('https://panamer.com', 'info@panamer.com')
    ,@('https://diagnosticoysoluciones.com', 'info@diagnosticoysoluciones.com')
    ,@('https://improntait.com', 'info@improntait.com')
    ,@('https://istea.com.ar', 'info@unlz.edu.ar, info@udelaciudad.edu.ar')
    ,@('https://swapconsulting.com.ar', 'info@swapconsulting.com.ar')
    ,@('https://avalith.net', 'hello@avalith.net')
    ,@('https://lanzallamas.tv', 'info@lanzallamas.tv')
    ,@('https://incutex.com.ar', 'info@incutex.com.ar')
    ,@('https://cactus.ws', 'info@cactus.ws, info@mysite.com')
    ,@('https://arenamaipu.com.ar', 'info@arenamaipu.com.ar')
    ,@('https://argeniss.com', 'info@argeniss.com')
)

$startRow = 193
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

